$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 11 (shifts old rows 11-16 down to rows 14-19).
$ws.Rows("11:13").Insert()

# Carry the column-A cell style (bold/centered/bordered) down onto the newly
# inserted rows by copying the formatted A10 cell into each of them.
$ws.Range("A10").Copy($ws.Range("A11"))
$ws.Range("A10").Copy($ws.Range("A12"))
$ws.Range("A10").Copy($ws.Range("A13"))

# --- Column A sequence numbers ---
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15
$ws.Range("A18").Value = 16
$ws.Range("A19").Value = 17

# --- Row labels (column B) ---
$ws.Range("B3").Value = "ND Single"
$ws.Range("B4").Value = "RD Single"
$ws.Range("B5").Value = "TD Single"
$ws.Range("B6").Value = "Morris"
$ws.Range("B7").Value = "Ring Perpendicular to ND"
$ws.Range("B8").Value = "Ring Perpendicular to RD"
$ws.Range("B9").Value = "Ring Perpendicular to TD"
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"

# --- Row 2 header labels (HKL indices) ---
$ws.Range("C2").Value = "[1, 1, 1]"
$ws.Range("D2").Value = "[2, 0, 0]"
$ws.Range("E2").Value = "[2, 2, 0]"
$ws.Range("F2").Value = "[3, 1, 1]"
$ws.Range("G2").Value = "[2, 2, 2]"
$ws.Range("H2").Value = "[4, 0, 0]"
$ws.Range("I2").Value = "[3, 3, 1]"
$ws.Range("J2").Value = "[4, 2, 0]"
$ws.Range("K2").Value = "[4, 2, 2]"
$ws.Range("L2").Value = "[5, 1, 1]"
$ws.Range("M2").Value = "[3, 3, 3]"
$ws.Range("N2").Value = "2Pairs"
$ws.Range("O2").Value = "4Pairs"
$ws.Range("P2").Value = "MaxUnique"

# --- Data values rows 10-19 (columns C-P) ---
# Row 10
$ws.Range("C10").Value = 1.038744308907424
$ws.Range("D10").Value = 0.9183273690739626
$ws.Range("E10").Value = 1.030356904684638
$ws.Range("F10").Value = 0.9715609114691947
$ws.Range("G10").Value = 1.038744308907424
$ws.Range("H10").Value = 0.9183273690739626
$ws.Range("I10").Value = 1.028838980602037
$ws.Range("J10").Value = 0.981524793711655
$ws.Range("K10").Value = 1.007886970994227
$ws.Range("L10").Value = 0.9385535419865516
$ws.Range("M10").Value = 1.038744308907424
$ws.Range("N10").Value = 0.9743421368793004
$ws.Range("O10").Value = 0.9897473735338049
$ws.Range("P10").Value = 0.9894742226787114

# Row 11
$ws.Range("C11").Value = 1.122270634628128
$ws.Range("D11").Value = 0.6037098106166887
$ws.Range("E11").Value = 1.174968010584092
$ws.Range("F11").Value = 0.9083047725391604
$ws.Range("G11").Value = 1.122270634628128
$ws.Range("H11").Value = 0.6037098106166887
$ws.Range("I11").Value = 1.149697570593438
$ws.Range("J11").Value = 0.9634065092288461
$ws.Range("K11").Value = 1.049001487012982
$ws.Range("L11").Value = 0.7471014611753369
$ws.Range("M11").Value = 1.122843606855655
$ws.Range("N11").Value = 0.8893389106003904
$ws.Range("O11").Value = 0.9523133070920174
$ws.Range("P11").Value = 0.964807532047334

# Row 12
$ws.Range("C12").Value = 1.122496269891132
$ws.Range("D12").Value = 0.6038607172078755
$ws.Range("E12").Value = 1.174645707861191
$ws.Range("F12").Value = 0.9083788971790877
$ws.Range("G12").Value = 1.122496269891132
$ws.Range("H12").Value = 0.6038607172078755
$ws.Range("I12").Value = 1.149610352880861
$ws.Range("J12").Value = 0.9633572317440933
$ws.Range("K12").Value = 1.049112769264316
$ws.Range("L12").Value = 0.747094450982202
$ws.Range("M12").Value = 1.12307343699413
$ws.Range("N12").Value = 0.8892532125345332
$ws.Range("O12").Value = 0.9523453980348217
$ws.Range("P12").Value = 0.9648195496263448

# Row 13
$ws.Range("C13").Value = 1.122137254422704
$ws.Range("D13").Value = 0.603695939250628
$ws.Range("E13").Value = 1.175071141235051
$ws.Range("F13").Value = 0.9083178956199833
$ws.Range("G13").Value = 1.122137254422704
$ws.Range("H13").Value = 0.603695939250628
$ws.Range("I13").Value = 1.149752908876161
$ws.Range("J13").Value = 0.963470783825203
$ws.Range("K13").Value = 1.048953763357867
$ws.Range("L13").Value = 0.7469951423105304
$ws.Range("M13").Value = 1.122714160319457
$ws.Range("N13").Value = 0.8893835402428394
$ws.Range("O13").Value = 0.9523055576320915
$ws.Range("P13").Value = 0.9647993536122659

# Row 14
$ws.Range("C14").Value = 1.214932
$ws.Range("D14").Value = 0.6329679999999995
$ws.Range("E14").Value = 1.083404
$ws.Range("F14").Value = 0.9198120000000004
$ws.Range("G14").Value = 1.214932
$ws.Range("H14").Value = 0.6329679999999995
$ws.Range("I14").Value = 1.113280000000003
$ws.Range("J14").Value = 0.9255720000000017
$ws.Range("K14").Value = 1.080879999999997
$ws.Range("L14").Value = 0.7792079999999997
$ws.Range("M14").Value = 1.215348000000002
$ws.Range("N14").Value = 0.8581859999999999
$ws.Range("O14").Value = 0.9627790000000001
$ws.Range("P14").Value = 0.9687570000000002

# Row 15
$ws.Range("C15").Value = 1.3
$ws.Range("D15").Value = 0.66
$ws.Range("E15").Value = 1
$ws.Range("F15").Value = 0.93
$ws.Range("G15").Value = 1.3
$ws.Range("H15").Value = 0.66
$ws.Range("I15").Value = 1.0809625
$ws.Range("J15").Value = 0.89
$ws.Range("K15").Value = 1.115974999999998
$ws.Range("L15").Value = 0.8093124999999993
$ws.Range("M15").Value = 1.3
$ws.Range("N15").Value = 0.8300000000000001
$ws.Range("O15").Value = 0.9725
$ws.Range("P15").Value = 0.9732812499999997

# Row 16
$ws.Range("C16").Value = 1.170613354495998
$ws.Range("D16").Value = 0.7996407009279987
$ws.Range("E16").Value = 0.9980243648512015
$ws.Range("F16").Value = 0.9572689466368016
$ws.Range("G16").Value = 1.170613354495998
$ws.Range("H16").Value = 0.7996407009279987
$ws.Range("I16").Value = 1.045232457113597
$ws.Range("J16").Value = 0.9341890803712021
$ws.Range("K16").Value = 1.064747163443193
$ws.Range("L16").Value = 0.8857089672192003
$ws.Range("M16").Value = 1.170624105267198
$ws.Range("N16").Value = 0.8988325328896001
$ws.Range("O16").Value = 0.9813868417280001
$ws.Range("P16").Value = 0.9819281293823992

# Row 17
$ws.Range("C17").Value = 0.9940391395698324
$ws.Range("D17").Value = 0.993949639068021
$ws.Range("E17").Value = 0.9943979510264693
$ws.Range("F17").Value = 0.993737914763454
$ws.Range("G17").Value = 0.9940391395698324
$ws.Range("H17").Value = 0.993949639068021
$ws.Range("I17").Value = 0.9939914371923674
$ws.Range("J17").Value = 0.9951644108813726
$ws.Range("K17").Value = 0.9940792992644015
$ws.Range("L17").Value = 0.993158591526912
$ws.Range("M17").Value = 0.9940818321432727
$ws.Range("N17").Value = 0.9941737950472451
$ws.Range("O17").Value = 0.9940311611069441
$ws.Range("P17").Value = 0.9940647979116037

# Row 18
$ws.Range("C18").Value = 0.9997851308671505
$ws.Range("D18").Value = 1.014956549687918
$ws.Range("E18").Value = 0.9782299330050362
$ws.Range("F18").Value = 0.9972592705738303
$ws.Range("G18").Value = 0.9997851308671505
$ws.Range("H18").Value = 1.014956549687918
$ws.Range("I18").Value = 0.9847985263926713
$ws.Range("J18").Value = 0.9918459484190866
$ws.Range("K18").Value = 0.9952113044551482
$ws.Range("L18").Value = 1.003641463572626
$ws.Range("M18").Value = 0.9997851308671505
$ws.Range("N18").Value = 0.9965932413464769
$ws.Range("O18").Value = 0.9975577210334836
$ws.Range("P18").Value = 0.9957160158716833

# Row 19
$ws.Range("C19").Value = 0.9666140924715699
$ws.Range("D19").Value = 1.030738206699461
$ws.Range("E19").Value = 0.9849797428625129
$ws.Range("F19").Value = 1.003058585190816
$ws.Range("G19").Value = 0.9666140924715699
$ws.Range("H19").Value = 1.030738206699461
$ws.Range("I19").Value = 0.979094158769794
$ws.Range("J19").Value = 1.002415833399573
$ws.Range("K19").Value = 0.9852438911320724
$ws.Range("L19").Value = 1.019044782625554
$ws.Range("M19").Value = 0.9666559414446625
$ws.Range("N19").Value = 1.007858974780987
$ws.Range("O19").Value = 0.9963476568060901
$ws.Range("P19").Value = 0.9963986616439193

Write-Host "Edit applied"